$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы"
$ws.Cells.Item(2, 2).Value = 81344
$ws.Cells.Item(3, 1).Value = "Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г"
$ws.Cells.Item(3, 2).Value = 2141
$ws.Cells.Item(4, 1).Value = "Пустырник трава 50г"
$ws.Cells.Item(4, 2).Value = 6454
$ws.Cells.Item(5, 1).Value = "Чабрец трава 50г"
$ws.Cells.Item(5, 2).Value = 12516
$ws.Cells.Item(6, 1).Value = "Багульник болотный побеги 50г"
$ws.Cells.Item(6, 2).Value = 8167
$ws.Cells.Item(7, 1).Value = "Валериана корневища с корнями 50г"
$ws.Cells.Item(7, 2).Value = 14250
$ws.Cells.Item(8, 1).Value = "Шалфей листья 50г"
$ws.Cells.Item(8, 2).Value = 26054
$ws.Cells.Item(9, 1).Value = "Дуба кора 75г"
$ws.Cells.Item(9, 2).Value = 52581
$ws.Cells.Item(10, 1).Value = "Бессмертник песчаный цветки 30г"
$ws.Cells.Item(10, 2).Value = 21241
$ws.Cells.Item(11, 1).Value = "Рябина плоды 50г"
$ws.Cells.Item(11, 2).Value = 1400
$ws.Cells.Item(12, 1).Value = "Крушина кора 50г"
$ws.Cells.Item(12, 2).Value = 8070
$ws.Cells.Item(13, 1).Value = "Лен семена 100г"
$ws.Cells.Item(13, 2).Value = 44676
$ws.Cells.Item(14, 1).Value = "Аир корневища 75г"
$ws.Cells.Item(14, 2).Value = 5631
$ws.Cells.Item(15, 1).Value = "Ромашка цветки вн 50г"
$ws.Cells.Item(15, 2).Value = 84528
$ws.Cells.Item(16, 1).Value = "Девясил корневища и корни 50г"
$ws.Cells.Item(16, 2).Value = 14247
$ws.Cells.Item(17, 1).Value = "Укроп пахучий плоды 50г"
$ws.Cells.Item(17, 2).Value = 52830
$ws.Cells.Item(18, 1).Value = "Сб. Грудной №4 50г"
$ws.Cells.Item(18, 2).Value = 29218
$ws.Cells.Item(19, 1).Value = "Мать-и-мачеха листья 35г"
$ws.Cells.Item(19, 2).Value = 22782
$ws.Cells.Item(20, 1).Value = "Ламинарии слоевища (морская капуста) 100г"
$ws.Cells.Item(20, 2).Value = 12806
$ws.Cells.Item(21, 1).Value = "Пижма цветки 75г"
$ws.Cells.Item(21, 2).Value = 13834
$ws.Cells.Item(22, 1).Value = "Полынь горькая трава 50г"
$ws.Cells.Item(22, 2).Value = 35936
$ws.Cells.Item(23, 1).Value = "Брусника листья 50г"
$ws.Cells.Item(23, 2).Value = 15271
$ws.Cells.Item(24, 1).Value = "Зверобой трава 50г"
$ws.Cells.Item(24, 2).Value = 33054
$ws.Cells.Item(25, 1).Value = "Череда трава 50г"
$ws.Cells.Item(25, 2).Value = 11421
$ws.Cells.Item(26, 1).Value = "Шиповник плоды низковитаминные 50г"
$ws.Cells.Item(26, 2).Value = 36352
$ws.Cells.Item(27, 1).Value = "Тысячелистник трава 50г"
$ws.Cells.Item(27, 2).Value = 15009
$ws.Cells.Item(28, 1).Value = "Можжевельник плоды 50г"
$ws.Cells.Item(28, 2).Value = 13788
$ws.Cells.Item(29, 1).Value = "Кукуруза столбики с рыльцами 40г"
$ws.Cells.Item(29, 2).Value = 28763
$ws.Cells.Item(30, 1).Value = "Береза почки 50г"
$ws.Cells.Item(30, 2).Value = 20157
$ws.Cells.Item(31, 1).Value = "Эвкалипт прутовидный листья 75г"
$ws.Cells.Item(31, 2).Value = 30343
$ws.Cells.Item(32, 1).Value = "Липа цветки 35г"
$ws.Cells.Item(32, 2).Value = 24702
$ws.Cells.Item(33, 1).Value = "Чага (березовый гриб) 50г"
$ws.Cells.Item(33, 2).Value = 33278
$ws.Cells.Item(34, 1).Value = "Сб. Фитонефрол (Урологический сбор) 50г"
$ws.Cells.Item(34, 2).Value = 21026
$ws.Cells.Item(35, 1).Value = "Спорыш трава 50г"
$ws.Cells.Item(35, 2).Value = 16268
$ws.Cells.Item(36, 1).Value = "Боярышник плоды 75г"
$ws.Cells.Item(36, 2).Value = 26958
$ws.Cells.Item(37, 1).Value = "Толокнянка листья 50г"
$ws.Cells.Item(37, 2).Value = 9480
$ws.Cells.Item(38, 1).Value = "Сенна листья 50г"
$ws.Cells.Item(38, 2).Value = 26459
$ws.Cells.Item(39, 1).Value = "Солодка корни 50г"
$ws.Cells.Item(39, 2).Value = 42981
$ws.Cells.Item(40, 1).Value = "Подорожник большой листья 50г"
$ws.Cells.Item(40, 2).Value = 12334
$ws.Cells.Item(41, 1).Value = "Эрва шерстистая трава 30г"
$ws.Cells.Item(41, 2).Value = 19219
$ws.Cells.Item(42, 1).Value = "Чистотел трава 50г"
$ws.Cells.Item(42, 2).Value = 24234
$ws.Cells.Item(43, 1).Value = "Крапива листья 50г"
$ws.Cells.Item(43, 2).Value = 18891
$ws.Cells.Item(44, 1).Value = "Алтей корни 75г"
$ws.Cells.Item(44, 2).Value = 8031
$ws.Cells.Item(45, 1).Value = "Мята перечная листья 50г"
$ws.Cells.Item(45, 2).Value = 33456
$ws.Cells.Item(46, 1).Value = "Сб. Фитопектол №1 (Грудной сбор №1) 35г"
$ws.Cells.Item(46, 2).Value = 7223
$ws.Cells.Item(47, 1).Value = "Ноготки цветки 50г"
$ws.Cells.Item(47, 2).Value = 36240
$ws.Cells.Item(48, 1).Value = "Сб. Фитопектол №2 (Грудной сбор №2) 35г"
$ws.Cells.Item(48, 2).Value = 9766
$ws.Cells.Item(49, 1).Value = "Фп Детский травяной чай `"ФармаЦветик® для иммунитета`" 20х1,5 г"
$ws.Cells.Item(50, 1).Value = "Фп Фиточай `"Лактафитол`" (БАД) 20х1,5 г"
$ws.Cells.Item(50, 2).Value = 14985
$ws.Cells.Item(51, 1).Value = "Фп Детский травяной чай `"ФармаЦветик®  при простуде`" 20х1,5 г"
$ws.Cells.Item(51, 2).Value = 3770
$ws.Cells.Item(52, 1).Value = "Фп Детский травяной чай `"ФармаЦветик® для животика`" 20х1,5 г"
$ws.Cells.Item(52, 2).Value = 5590
$ws.Cells.Item(53, 1).Value = "Фп Детский травяной чай `"ФармаЦветик® для спокойного сна`" 20х1,5 г"
$ws.Cells.Item(53, 2).Value = 8468
$ws.Cells.Item(54, 1).Value = "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем`"(БАД) 20*1,5г"
$ws.Cells.Item(54, 2).Value = 7700
$ws.Cells.Item(55, 1).Value = "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем`" (БАД) 20*1,5г"
$ws.Cells.Item(55, 2).Value = 9100
$ws.Cells.Item(56, 1).Value = "Фп `"Щедрость природы® Фиточай для иммунитета`" 20х2,0 г"
$ws.Cells.Item(56, 2).Value = 432
$ws.Cells.Item(57, 1).Value = "Фп `"Щедрость природы® Фиточай при простуде`" 20х2,0 г"
$ws.Cells.Item(57, 2).Value = 594
$ws.Cells.Item(58, 1).Value = "Фп `"Щедрость природы® Фиточай кардиологический`" 20х2,0 г"
$ws.Cells.Item(58, 2).Value = 1026
$ws.Cells.Item(59, 1).Value = "Фп `"Щедрость природы® Фиточай успокоительный`"20х2,0 г"
$ws.Cells.Item(59, 2).Value = 1350
$ws.Cells.Item(60, 1).Value = "Фп `"Щедрость природы® Фиточай диабетический`" 20х2,0 г"
$ws.Cells.Item(60, 2).Value = 828
$ws.Cells.Item(61, 1).Value = "Фп Шалфей листья 20х1,5г"
$ws.Cells.Item(61, 2).Value = 108876
$ws.Cells.Item(62, 1).Value = "Фп Сб. Грудной №4 20x2,0г"
$ws.Cells.Item(62, 2).Value = 509468
$ws.Cells.Item(63, 1).Value = "Фп `"Щедрость природы® Фиточай очищающий`" 20х2,0 г"
$ws.Cells.Item(63, 2).Value = 1494
$ws.Cells.Item(64, 1).Value = "Фп Фиточай `"Опалиховский`" (БАД) 20х2,0 г"
$ws.Cells.Item(64, 2).Value = 4572
$ws.Cells.Item(65, 1).Value = "Фп Фиточай `"Тибетский`" (БАД) 20х2,0  г"
$ws.Cells.Item(65, 2).Value = 9395
$ws.Cells.Item(66, 1).Value = "Фп `"Щедрость природы® Фиточай для пищеварения`" 20х2,0 г"
$ws.Cells.Item(66, 2).Value = 1602
$ws.Cells.Item(67, 1).Value = "Фп Чистотел трава 20х1,5г"
$ws.Cells.Item(67, 2).Value = 24684
$ws.Cells.Item(68, 1).Value = "Фп Подорожник листья 20x1,5г"
$ws.Cells.Item(68, 2).Value = 24506
$ws.Cells.Item(69, 1).Value = "Фп Пустырник трава 20x1,5г"
$ws.Cells.Item(69, 2).Value = 35688
$ws.Cells.Item(70, 1).Value = "Фп Брусника листья 20х1,5г"
$ws.Cells.Item(70, 2).Value = 70019
$ws.Cells.Item(71, 1).Value = "Фп Сб. Бруснивер 20x2,0г"
$ws.Cells.Item(71, 2).Value = 189546
$ws.Cells.Item(72, 1).Value = "Фп Пастушья сумка трава 20х1,5г"
$ws.Cells.Item(72, 2).Value = 4822
$ws.Cells.Item(73, 1).Value = "Фп Ромашка цветки 20x1,5г"
$ws.Cells.Item(73, 2).Value = 1299546
$ws.Cells.Item(74, 1).Value = "Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г"
$ws.Cells.Item(74, 2).Value = 172471
$ws.Cells.Item(75, 1).Value = "Фп Сб. Элекасол 20x2,0г"
$ws.Cells.Item(75, 2).Value = 37854
$ws.Cells.Item(76, 1).Value = "Фп Череда трава 20х1,5г"
$ws.Cells.Item(76, 2).Value = 47897
$ws.Cells.Item(77, 1).Value = "Фп Сенна листья 20x1,5г"
$ws.Cells.Item(77, 2).Value = 71090
$ws.Cells.Item(78, 1).Value = "Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г"
$ws.Cells.Item(78, 2).Value = 22822
$ws.Cells.Item(79, 1).Value = "Фп Мелисса лекарственная трава 20x1,5г"
$ws.Cells.Item(79, 2).Value = 37224
$ws.Cells.Item(80, 1).Value = "Фп Шиповник плоды 20х2,0г"
$ws.Cells.Item(80, 2).Value = 50346
$ws.Cells.Item(81, 1).Value = "Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г"
$ws.Cells.Item(81, 2).Value = 89595
$ws.Cells.Item(82, 1).Value = "Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г"
$ws.Cells.Item(82, 2).Value = 92551
$ws.Cells.Item(83, 1).Value = "Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г"
$ws.Cells.Item(83, 2).Value = 89979
$ws.Cells.Item(84, 1).Value = "Фп Толокнянка листья 20x1,5г"
$ws.Cells.Item(84, 2).Value = 42208
$ws.Cells.Item(85, 1).Value = "Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г"
$ws.Cells.Item(85, 2).Value = 58159
$ws.Cells.Item(86, 1).Value = "Фп Крапива листья 20x1,5г"
$ws.Cells.Item(86, 2).Value = 69149
$ws.Cells.Item(87, 1).Value = "Фп Зверобой трава 20x1,5г"
$ws.Cells.Item(87, 2).Value = 54293
$ws.Cells.Item(88, 1).Value = "Фп Чабрец трава 20x1,5 г"
$ws.Cells.Item(88, 2).Value = 71604
$ws.Cells.Item(89, 1).Value = "Фп Душица трава 20x1,5г"
$ws.Cells.Item(89, 2).Value = 30582
$ws.Cells.Item(90, 1).Value = "Фп Хвощ полевой трава 20х1,5г"
$ws.Cells.Item(90, 2).Value = 30454
$ws.Cells.Item(91, 1).Value = "Фп Сб. Желудочный №3 20x2,0г"
$ws.Cells.Item(91, 2).Value = 25902
$ws.Cells.Item(92, 1).Value = "Фп Мята перечная листья 20x1,5г"
$ws.Cells.Item(92, 2).Value = 77945
$ws.Cells.Item(93, 1).Value = "Фп Липа цветки 20x1,5г"
$ws.Cells.Item(93, 2).Value = 82023
$ws.Cells.Item(94, 1).Value = "Фп Сб. Арфазетин-Э 20x2,0г"
$ws.Cells.Item(94, 2).Value = 58481
$ws.Cells.Item(95, 1).Value = "Фп Береза листья 20x1,5г"
$ws.Cells.Item(95, 2).Value = 5926
$ws.Cells.Item(96, 1).Value = "Фп Золототысячник трава 20х1,5г"
$ws.Cells.Item(96, 2).Value = 5475
$ws.Cells.Item(97, 1).Value = "Фп Фиалка трехцветная трава 20x1,5г"
$ws.Cells.Item(97, 2).Value = 5148
$ws.Cells.Item(98, 1).Value = "Фп Боярышник плоды 20х3,0г"
$ws.Cells.Item(98, 2).Value = 25460
$ws.Cells.Item(99, 1).Value = "Фп Аир корневища 20x1,5г"
$ws.Cells.Item(99, 2).Value = 6323
$ws.Cells.Item(100, 1).Value = "Фп Пижма цветки 20х1,5г"
$ws.Cells.Item(100, 2).Value = 11118
$ws.Cells.Item(101, 1).Value = "Фп Ольха соплодия 20х1,5г"
$ws.Cells.Item(101, 2).Value = 5839
$ws.Cells.Item(102, 1).Value = "Фп Ноготки цветки 20x1,5г"
$ws.Cells.Item(102, 2).Value = 26537
$ws.Cells.Item(103, 1).Value = "Фп Девясил корневища и корни 20х1,5г"
$ws.Cells.Item(103, 2).Value = 10134
$ws.Cells.Item(104, 1).Value = "Фп Дуб кора 20х1,5г"
$ws.Cells.Item(104, 2).Value = 7245
$ws.Cells.Item(105, 1).Value = "Фп Крушина кора 20x1,5г"
$ws.Cells.Item(105, 2).Value = 10964
$ws.Cells.Item(106, 1).Value = "Фп Бадан корневища 20x1,5г"
$ws.Cells.Item(106, 2).Value = 1627
$ws.Cells.Item(107, 1).Value = "Фп Валериана корневища с корнями 20x1,5г"
$ws.Cells.Item(107, 2).Value = 29960
$ws.Cells.Item(108, 1).Value = "Фп Тысячелистник трава 20x1,5г"
$ws.Cells.Item(108, 2).Value = 28130
$ws.Cells.Item(109, 1).Value = "Фп Кровохлебка корневища и корни 20x1,5г"
$ws.Cells.Item(109, 2).Value = 10042
$ws.Cells.Item(110, 1).Value = "Фп Лапчатка корневища 20x2,5г"
$ws.Cells.Item(110, 2).Value = 6421
$ws.Cells.Item(111, 1).Value = "Фп Почечный чай листья 20x1,5г"
$ws.Cells.Item(111, 2).Value = 148743

$ws.Range("A16").Select()
